$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions scheduled update)

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.232.98'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.65%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.676.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.49%  '

# Row 4
$ws.Range('E4').Value = '  +0.34%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.25%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5327'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.38%  '

# Row 7
$ws.Range('E7').Value = '  +0.30%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2676'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.19%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06473'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.05%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.87'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.22%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07525'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.41%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.672.05'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.62%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.513'
$ws.Range('D13').Style = 'Normal'

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5758'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.44%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008455'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.73%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.60'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.52%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.290.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.71%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.898'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.84%  '

# Row 19
$ws.Range('E19').Value = '  +0.28%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.92%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.08%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.191'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.54%  '

# Row 23
$ws.Range('E23').Value = '  +0.36%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.02%  '

# Row 25
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.820'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.91%  '

# Row 26
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1272'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.24%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.87%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06438'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.29%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.379'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.51%  '

# Row 30
$ws.Range('E30').Value = '  +0.19%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.582'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.59%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.591'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.07%  '

# Row 33
$ws.Range('E33').Value = '  +0.74%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.030'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.11%  '

# Row 35
$ws.Range('E35').Value = '  +1.49%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.403'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.52%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.724'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.57%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.255'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.41%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.112.51'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.37%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01621'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.19%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8725'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.03%  '

# Row 42
$ws.Range('E42').Value = '  +0.50%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.37'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.18%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.827.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.65%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000109'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.26%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.09%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.155'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.99%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.55%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05258'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.78%  '

# Row 50
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.091'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.98%  '

# Row 51
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4289'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.13%  '
